$wb = $excel.ActiveWorkbook

function Add-DataRow($ws, $row, $timeVal, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    # Column A: date/time serial, formatted like the row above it
    $ws.Cells.Item($row, 1).Value = $timeVal
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

    # Columns B-E: hex byte-string text fields
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal

    # Column F: plain integer
    $ws.Cells.Item($row, 6).Value = $fVal

    # Column G: numeric, except when the caller passes an oversize digit
    # string that must be preserved verbatim as text (beyond double
    # precision), matching how the source data stores it.
    if ($gVal -is [string]) {
        $gCell = $ws.Cells.Item($row, 7)
        $gCell.NumberFormat = "@"
        $gCell.Value = $gVal
        $gCell.Style = "Normal"
    } else {
        $ws.Cells.Item($row, 7).Value = $gVal
    }

    # Columns H-I: plain integers
    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

$bigG = [double]"5.68631262647114e+23"

$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-DataRow $ws1 51 45749.69065453704 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x62" "0xe" 400 $bigG 354 14

$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-DataRow $ws2 53 45749.66216435185 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x6a" "0x19" 400 "568631262647113771663628" 362 25

$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-DataRow $ws3 51 45749.72177849537 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x62" "0x14" 400 $bigG 354 20

$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-DataRow $ws4 51 45749.85645540509 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x6a" "0x19" 400 $bigG 362 25
